# Daily "water delivery countdown" update:
# For every data row (2..last used row), column E holds the remaining
# count and column F holds the cycle start date (yyyyMMdd, numeric).
# Each run decrements E by 1. If that would bring E to 0, the cycle
# restarts: E resets to the row's total (column D) and F is pushed
# forward by D days (a new delivery cycle begins).
# Rows whose F value isn't a well-formed 8-digit date (e.g. a typo)
# are left untouched, since the countdown can't be computed for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $total = $ws.Cells.Item($row, 4).Value2      # column D - 总天 (total days)
    $remain = $ws.Cells.Item($row, 5).Value2      # column E - 剩余 (remaining)
    $start = $ws.Cells.Item($row, 6).Value2      # column F - 开始时间 (start date)

    if ($null -eq $total -or $null -eq $remain -or $null -eq $start) {
        continue
    }

    $startStr = [string]([int64]$start)
    if ($startStr.Length -ne 8) {
        # malformed date, skip this row
        continue
    }

    $newRemain = $remain - 1

    if ($newRemain -le 0) {
        $year = [int]$startStr.Substring(0, 4)
        $month = [int]$startStr.Substring(4, 2)
        $day = [int]$startStr.Substring(6, 2)
        $startDate = Get-Date -Year $year -Month $month -Day $day
        $newStartDate = $startDate.AddDays([double]$total)
        $newStart = [int64]$newStartDate.ToString("yyyyMMdd")

        $ws.Cells.Item($row, 5).Value = $total
        $ws.Cells.Item($row, 6).Value = $newStart
    }
    else {
        $ws.Cells.Item($row, 5).Value = $newRemain
    }
}
